# Adding test case to authoring
# Adds a new authoring test case row (VerifyDraftPostDisplayInUserOwnProfile /
# OPQA-1197) to the "Test Cases" sheet, harmonises a handful of cell styles
# in that area of the sheet (hyperlink borders, Runmode-column formatting),
# flips several Results cells from PASS to SKIP, and updates the sheet's
# selection to the newly active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# ---------------------------------------------------------------------
# 1) Harmonise the hyperlink (JIRA ID) cell formatting in column B so the
#    bordered Hyperlink style already used by the later rows (e.g. B50)
#    is also used by these earlier rows.
# ---------------------------------------------------------------------
$ws.Range("B50").Copy()
$hyperlinkCells = @("B8","B12","B13","B15","B16","B18","B19","B20","B21","B22","B28","B34","B35","B36","B37","B38","B39")
foreach ($addr in $hyperlinkCells) {
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Harmonise the Runmode column (D) formatting for rows 59-63 so they
#    match the rest of the column (D2:D58 already use this style).
# ---------------------------------------------------------------------
$ws.Range("D58").Copy()
$runmodeCells = @("D59","D60","D61","D62","D63")
foreach ($addr in $runmodeCells) {
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) Update several Results cells (column E) from PASS to SKIP.
# ---------------------------------------------------------------------
$skipCells = @("E39","E50","E51","E53","E54","E55","E56","E57","E58","E59","E61","E62","E63")
foreach ($addr in $skipCells) {
    $ws.Range($addr).Value = "SKIP"
}

# ---------------------------------------------------------------------
# 4) Append the new authoring test case as row 64. Formats are copied
#    in first (bordered cell style, matching the rest of the table; the
#    Description cell also gets wrap-text like the rest of column C),
#    then the values are written in Description -> JIRA ID -> TCID
#    order so any newly-created shared strings line up the same way
#    the authoring tool produced them.
# ---------------------------------------------------------------------
$ws.Range("A59").Copy()
$newRowCells = @("A64","B64","D64","E64")
foreach ($addr in $newRowCells) {
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

$ws2 = $wb.Worksheets.Item("Test Case Steps")
$ws2.Range("D2").Copy()
$ws.Range("C64").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("C64").Value = "Verify that Draft Post tab is displayed only in the users own profile and only `nwhen the user has at least one draft post"
$ws.Range("B64").Value = "OPQA-1197"
$ws.Range("A64").Value = "VerifyDraftPostDisplayInUserOwnProfile"
$ws.Range("D64").Value = "Y"
$ws.Range("E64").Value = "PASS"

$ws.Rows.Item(64).RowHeight = 30

# ---------------------------------------------------------------------
# 5) Update the sheet selection to reflect the newly active cell.
# ---------------------------------------------------------------------
$ws.Range("D61").Select()
